# Update the "取得日時" (retrieved timestamp) column on the active sheet
# ("ランサーズ" / 案件情報) for all existing data rows from the old
# timestamp "2025-09-12 12:34:47" to the new one "2025-09-12 12:43:42".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$oldValue = "2025-09-12 12:34:47"
$newValue = "2025-09-12 12:43:42"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value = $newValue
    }
}
